# Automatische test-sync: 2025-06-24 20:12:50
# Appends a new mail-log entry (row 19) to the "Logs" sheet and updates the
# corresponding category count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 19

$logs.Cells.Item($newRow, 1).Value = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Mijn bestelling is incompleet geleverd. Graag hoor ik hoe dit wordt opgelost."
$logs.Cells.Item($newRow, 4).Value = "Bestelling / Levering"
$logs.Cells.Item($newRow, 5).Value = "Beste klant,`nBedankt voor uw bericht. Het spijt me te horen dat uw bestelling incompleet geleverd is. Om dit probleem op te lossen, hebben we wat meer informatie nodig. Kunt u ons uw bestelnummer doorgeven, zodat we kunnen controleren wat er precies misgaan is? Op basis daarvan zullen we zo snel mogelijk een oplossing voor u vinden.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Cells.Item($newRow, 6).Value = "2025-06-24 20:12:39"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# The multi-line text in column E triggers an automatic row-height override;
# AutoFit restores the row to the sheet's default (unset) height, matching
# the plain <row r="19"> (no explicit ht/customHeight) seen in the other rows.
$logs.Rows.Item($newRow).AutoFit()

# Extend the conditional formatting ranges so the new row is covered too.
foreach ($fc in $logs.Range("D2:D18").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("D2:D19"))
}
foreach ($fc in $logs.Range("G2:G18").FormatConditions) {
    $fc.ModifyAppliesToRange($logs.Range("G2:G19"))
}

# Update the Dashboard count for "Bestelling / Levering" (row 8, column B).
$dashboard.Cells.Item(8, 2).Value = 2
